# Apply the commit's change: the filter breed used for the count-query
# ("StatOutput"/"StatOutput_Message" sheets) is switched from a stale
# breed ("Akita") to "Irish Setter", and the corresponding summary
# counts on the StatOutput sheet are updated to "0".

$wb = $excel.ActiveWorkbook

# --- StatOutput!A2 / StatOutput!B2 : "1"/"2" -> "0"/"0" -------------------
$wsStat = $wb.Worksheets.Item("StatOutput")

# Force text storage (matches the source file, where these are shared
# strings rather than numbers) by marking the cells as Text before
# writing the numeric-looking value.
$wsStat.Range("A2").NumberFormat = "@"
$wsStat.Range("A2").Value = "0"

$wsStat.Range("B2").NumberFormat = "@"
$wsStat.Range("B2").Value = "0"

# --- StatOutput_Message!A18 : Akita count-query -> Irish Setter count-query
$wsMsg = $wb.Worksheets.Item("StatOutput_Message")
$wsMsg.Range("A18").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Irish Setter']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
